$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-30 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-01 Tuesday", 2)
$d.Content.Find.Execute("507÷4=126, 3", $true, $false, $false, $false, $false, $true, 1, $false, "552÷6=92, 0", 2)
$d.Content.Find.Execute("397÷9=44, 1", $true, $false, $false, $false, $false, $true, 1, $false, "633÷9=70, 3", 2)
$d.Content.Find.Execute("966÷4=241, 2", $true, $false, $false, $false, $false, $true, 1, $false, "844÷4=211, 0", 2)
$d.Content.Find.Execute("224÷2=112, 0", $true, $false, $false, $false, $false, $true, 1, $false, "394÷3=131, 1", 2)
$d.Content.Find.Execute("218÷2=109, 0", $true, $false, $false, $false, $false, $true, 1, $false, "118÷8=14, 6", 2)
$d.Content.Find.Execute("544÷3=181, 1", $true, $false, $false, $false, $false, $true, 1, $false, "735÷7=105, 0", 2)
$d.Content.Find.Execute("793÷4=198, 1", $true, $false, $false, $false, $false, $true, 1, $false, "219÷9=24, 3", 2)
$d.Content.Find.Execute("316÷6=52, 4", $true, $false, $false, $false, $false, $true, 1, $false, "507÷7=72, 3", 2)
$d.Content.Find.Execute("359÷5=71, 4", $true, $false, $false, $false, $false, $true, 1, $false, "581÷3=193, 2", 2)
$d.Content.Find.Execute("992÷9=110, 2", $true, $false, $false, $false, $false, $true, 1, $false, "698÷2=349, 0", 2)
$d.Content.Find.Execute("659÷7=94, 1", $true, $false, $false, $false, $false, $true, 1, $false, "413÷3=137, 2", 2)
$d.Content.Find.Execute("143÷8=17, 7", $true, $false, $false, $false, $false, $true, 1, $false, "282÷7=40, 2", 2)
$d.Content.Find.Execute("134÷6=22, 2", $true, $false, $false, $false, $false, $true, 1, $false, "621÷4=155, 1", 2)
$d.Content.Find.Execute("606÷3=202, 0", $true, $false, $false, $false, $false, $true, 1, $false, "664÷2=332, 0", 2)
$d.Content.Find.Execute("262÷2=131, 0", $true, $false, $false, $false, $false, $true, 1, $false, "611÷6=101, 5", 2)
$d.Content.Find.Execute("221÷2=110, 1", $true, $false, $false, $false, $false, $true, 1, $false, "350÷3=116, 2", 2)
$d.Content.Find.Execute("433÷6=72, 1", $true, $false, $false, $false, $false, $true, 1, $false, "458÷4=114, 2", 2)
$d.Content.Find.Execute("548÷5=109, 3", $true, $false, $false, $false, $false, $true, 1, $false, "551÷3=183, 2", 2)
$d.Content.Find.Execute("558÷4=139, 2", $true, $false, $false, $false, $false, $true, 1, $false, "732÷6=122, 0", 2)
$d.Content.Find.Execute("890÷9=98, 8", $true, $false, $false, $false, $false, $true, 1, $false, "765÷2=382, 1", 2)
$d.Content.Find.Execute("207÷7=29, 4", $true, $false, $false, $false, $false, $true, 1, $false, "839÷2=419, 1", 2)
$d.Content.Find.Execute("916÷3=305, 1", $true, $false, $false, $false, $false, $true, 1, $false, "808÷7=115, 3", 2)
$d.Content.Find.Execute("795÷5=159, 0", $true, $false, $false, $false, $false, $true, 1, $false, "812÷6=135, 2", 2)
$d.Content.Find.Execute("127÷3=42, 1", $true, $false, $false, $false, $false, $true, 1, $false, "903÷9=100, 3", 2)
$d.Content.Find.Execute("480÷2=240, 0", $true, $false, $false, $false, $false, $true, 1, $false, "441÷2=220, 1", 2)
